$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Combine Name,Age into column A as a single CSV-style string per row (like exporting to CSV),
# then remove the now-redundant column B.
$rows = $ws.UsedRange.Rows.Count
for ($r = 1; $r -le $rows; $r++) {
    $nameVal = $ws.Cells.Item($r, 1).Value2
    $ageVal = $ws.Cells.Item($r, 2).Value2
    $ws.Cells.Item($r, 1).Value = "$nameVal,$ageVal"
}

$ws.Columns.Item(2).Delete()

$ws.Range("H2").Select()
